$wb = $excel.ActiveWorkbook

function Set-TextCell($sheet, $row, $col, $val) {
    # Leading apostrophe forces Excel to store the value as text (keeps
    # leading zeros / decimal-looking strings like "5.05" as strings,
    # matching the source data which uses t="inlineStr" for these columns).
    $sheet.Cells.Item($row, $col).Value = "'" + $val
    # Drop the quote-prefix formatting flag that the apostrophe trick adds
    # so the cell ends up with the workbook's default (no) style, same as
    # the target.
    $sheet.Cells.Item($row, $col).ClearFormats()
}

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet, positioned right after "总计"
#    and before "2022-Q3" (mirrors the workbook.xml sheet-order change).
#    Duplicating "2022-Q3" (instead of Worksheets.Add) keeps all of its
#    cell/row formatting (bold header, borders, index-column style)
#    intact, which a blank Add() + manual styling could not reliably
#    reproduce cross-sheet.
# ---------------------------------------------------------------------
$sheetQ3 = $wb.Worksheets.Item("2022-Q3")
$sheetQ3.Copy($sheetQ3)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# Row 2 already carries the right formatting (copied from 2022-Q3's row
# 2); just overwrite its values with the first 2022-Q4 holding.
$newSheet.Cells.Item(2,1).Value = 0
Set-TextCell $newSheet 2 2 "011410"
Set-TextCell $newSheet 2 3 "中信建投量化进取6个月持有期混合A"
Set-TextCell $newSheet 2 4 "5.05"
Set-TextCell $newSheet 2 5 "69.90"
Set-TextCell $newSheet 2 6 "0.43"
Set-TextCell $newSheet 2 7 "0.0217"
$newSheet.Cells.Item(2,8).Value = 9

# Add three more data rows (3-5), inheriting column A's style from row 2.
$newSheet.Rows.Item(3).Insert()
$newSheet.Rows.Item(4).Insert()
$newSheet.Rows.Item(5).Insert()

$newSheet.Cells.Item(2,1).Copy()
$newSheet.Range("A3:A5").PasteSpecial(-4122)

# Row 3
$newSheet.Cells.Item(3,1).Value = 1
Set-TextCell $newSheet 3 2 "013242"
Set-TextCell $newSheet 3 3 "北信瑞丰优势行业股票"
Set-TextCell $newSheet 3 4 "0.49"
Set-TextCell $newSheet 3 5 "91.56"
Set-TextCell $newSheet 3 6 "1.85"
Set-TextCell $newSheet 3 7 "0.0091"
$newSheet.Cells.Item(3,8).Value = 3

# Row 4
$newSheet.Cells.Item(4,1).Value = 2
Set-TextCell $newSheet 4 2 "004730"
Set-TextCell $newSheet 4 3 "建信量化事件驱动股票"
Set-TextCell $newSheet 4 4 "0.46"
Set-TextCell $newSheet 4 5 "90.94"
Set-TextCell $newSheet 4 6 "1.83"
Set-TextCell $newSheet 4 7 "0.0084"
$newSheet.Cells.Item(4,8).Value = 9

# Row 5
$newSheet.Cells.Item(5,1).Value = 3
Set-TextCell $newSheet 5 2 "011411"
Set-TextCell $newSheet 5 3 "中信建投量化进取6个月持有期混合C"
Set-TextCell $newSheet 5 4 "1.56"
Set-TextCell $newSheet 5 5 "69.90"
Set-TextCell $newSheet 5 6 "0.43"
Set-TextCell $newSheet 5 7 "0.0067"
$newSheet.Cells.Item(5,8).Value = 9

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row 2 for the
#    2022-Q4 totals, shifting the existing quarters down by one row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# The inserted row inherits row-1 (header) formatting from Excel's
# auto-fill-down behaviour; strip it so B2:D2 end up with the default
# (no) style, matching the source data rows.
$total.Range("B2:D2").ClearFormats()

# Give the new row's index cell (A2) the same style as the rest of the
# index column (copy from what is now A3, the old A2 / "2022-Q3" row).
$total.Cells.Item(3,1).Copy()
$total.Cells.Item(2,1).PasteSpecial(-4122)

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 4
$total.Cells.Item(2,4).Value = 0.05

# The original index column values (A3:A8) are left untouched by the
# insert; only the last row's ("2020-Q4", pushed from row 8 to row 9)
# index bumps from 6 to 7 in the target data.
$total.Cells.Item(9,1).Value = 7

# Restore the original active sheet/selection (sheet-copy/insert above
# leaves the newly created sheet selected).
$total.Activate()
[void]$total.Range("A1").Select()

Write-Output "2022-Q4 sheet + summary row added"
